$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list
# with the latest scraped figures. Values that would otherwise be
# mis-parsed by Excel as numbers (losing trailing zeros / punctuation)
# are written with a leading apostrophe to force plain text, matching
# the original inline-string cell content.
$ws.Range("D2").Value = "62.016.97"
$ws.Range("D3").Value = "2.421.57"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'562.89"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "'143.19"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.531"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").Value = "2.235.70"
$ws.Range("E9").Value = "  -7.55%  "
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("E12").Value = "  -4.00%  "
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").Value = "'26.15"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "'0.0000174"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "61.938.66"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "2.427.60"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "'11.31"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").Value = "'323.92"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Value = "'6.82"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").Value = "'4.13"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'67.05"
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("D25").Value = "'1.74"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").Value = "'8.77"
$ws.Range("E26").Value = "  -2.62%  "
$ws.Range("D27").Value = "'555.19"
$ws.Range("E27").Value = "  -5.40%  "
$ws.Range("D28").Value = "2.543.34"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "0.0₃0930"
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("D31").Value = "'8.20"
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("E32").Value = "  -4.67%  "
$ws.Range("E33").Value = "  -2.09%  "
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("D35").Value = "'1.50"
$ws.Range("E35").Value = "  -3.65%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "'4.74"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").Value = "'0.378"
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("D39").Value = "'5.45"
$ws.Range("E39").Value = "  -4.86%  "
$ws.Range("D40").Value = "'152.69"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").Value = "'18.64"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").Value = "'1.81"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").Value = "'0.993"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "'147.36"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("E45").Value = "  -5.00%  "
$ws.Range("D46").Value = "'3.63"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").Value = "'0.0528"
$ws.Range("E47").Value = "  -2.18%  "
$ws.Range("D48").Value = "'0.595"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").Value = "'19.79"
$ws.Range("E49").Value = "  -3.00%  "
$ws.Range("D50").Value = "'0.0919"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("E51").Value = "  -0.46%  "
